# Update rules in DiscountRules.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the rule table with a new (empty, formatted) column E, matching
# the existing column D formatting across the whole used range.
$ws.Range("D1:D26").Copy()
$ws.Range("E1:E26").PasteSpecial(-4122)

# Mirror the ACTION / Test values into the new column for the header rows
$ws.Range("E18").Value = $ws.Range("D18").Value()
$ws.Range("E19").Value = $ws.Range("C19").Value()

# New rule-row values in column C
$ws.Range("C25").Value = ""
$ws.Range("C26").Value = "TEST GIT SAMPLE"
